# reviewdb.xlsx: keep only the header row and the last review (helix jump /
# snizzvered@gmail.com / krigelron@gmail.com), dropping the first three
# review rows (block.chain.technology + com.singleton.helix x2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the formatting of C5/D5 (the row that survives the deletion below and
# becomes the new row 2) in far-away scratch cells so it can be restored
# later -- Hyperlinks.Add() (used further down) stamps its own generic
# "Hyperlink" style on the target cells, which we don't want here since the
# original workbook already carried bespoke styling on those columns.
$ws.Range("C5").Copy()
$ws.Range("Z100").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D5").Copy()
$ws.Range("Z101").PasteSpecial(-4122)   # xlPasteFormats

# Delete the first three review rows (2, 3, 4). Row 5 shifts up and becomes
# the new row 2, row 1 (the header) is untouched.
$ws.Rows("2:4").Delete()

# The worksheet's Hyperlinks collection still references the old row
# positions (C2:D2, C3:D3, C4:D4, C5:D5) and doesn't auto-adjust after the
# row delete, so rebuild it from scratch with just the two links that are
# still meaningful, now living on row 2.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:snizzvered@gmail.com", [Type]::Missing, [Type]::Missing, "snizzvered@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:krigelron@gmail.com", [Type]::Missing, [Type]::Missing, "krigelron@gmail.com")

# Restore the original cell formatting on C2/D2 that Hyperlinks.Add()
# overwrote. The scratch cells shifted up by 3 rows along with everything
# else from the earlier row delete, so they now live at Z97/Z98.
$ws.Range("Z97").Copy()
$ws.Range("C2").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("Z98").Copy()
$ws.Range("D2").PasteSpecial(-4122)     # xlPasteFormats

# Clean up the scratch cells so they don't leak into the saved sheet.
$ws.Range("Z97:Z98").Clear()

# Match the author's final selection (cell A2).
$ws.Range("A2").Select() | Out-Null
